$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()),
                (''selector'', RandomUnderSampler(random_state=42)),
                (''model'',
                 AdaBoostClassifier(estimator=DecisionTreeClassifier(criterion=''entropy'',
                                                                     max_depth=5,
                                                                     max_features=''sqrt'',
                                                                     random_state=42),
                                    n_estimators=10, random_state=42))])'
$ws.Range("B2").Value = 0.6319580419580418
$ws.Range("C2").Value = '{''selector'': RandomUnderSampler(random_state=42), ''scaler'': MinMaxScaler(), ''model__n_estimators'': 10, ''model__estimator__min_samples_split'': 2, ''model__estimator__min_samples_leaf'': 1, ''model__estimator__max_features'': ''sqrt'', ''model__estimator__max_depth'': 5, ''model__estimator__criterion'': ''entropy'', ''model__estimator__class_weight'': None}'
$ws.Range("D2").Value = 0.9237305050552113
$ws.Range("E2").Value = 0.4737500444000445
$ws.Range("F2").Value = 0.7272727272727272
$ws.Range("G2").Value = 0.9881859095293996
$ws.Range("H2").Value = 0.5311805555555555
$ws.Range("I2").Value = 0.7058823529411765
$ws.Range("J2").Value = 0.8685106382978722
$ws.Range("K2").Value = 0.4490000000000001
$ws.Range("L2").Value = 0.75
$ws.Range("M2").Value = '[1 0 1 1 1 1 0 1 0 1 0 1 0 1 1 0 0 1 1 1 1 0 1 1]'
$ws.Range("N2").Value = '[1 0 1 1 1 1 1 0 1 1 1 1 1 0 1 0 0 1 0 1 1 1 1 0]'

$ws.Range("A3").Value = 'Pipeline(steps=[(''scaler'', StandardScaler()),
                (''selector'', RandomUnderSampler(random_state=42)),
                (''model'',
                 AdaBoostClassifier(estimator=DecisionTreeClassifier(class_weight=''balanced'',
                                                                     max_depth=1,
                                                                     max_features=''sqrt'',
                                                                     min_samples_leaf=5,
                                                                     random_state=42),
                                    random_state=42))])'
$ws.Range("B3").Value = 0.6254778554778554
$ws.Range("C3").Value = '{''selector'': RandomUnderSampler(random_state=42), ''scaler'': StandardScaler(), ''model__n_estimators'': 50, ''model__estimator__min_samples_split'': 2, ''model__estimator__min_samples_leaf'': 5, ''model__estimator__max_features'': ''sqrt'', ''model__estimator__max_depth'': 1, ''model__estimator__criterion'': ''gini'', ''model__estimator__class_weight'': ''balanced''}'
$ws.Range("D3").Value = 0.9149271946510816
$ws.Range("E3").Value = 0.4634015040515041
$ws.Range("F3").Value = 0.6206896551724138
$ws.Range("G3").Value = 0.9842006407077476
$ws.Range("H3").Value = 0.5814345238095239
$ws.Range("I3").Value = 0.6923076923076923
$ws.Range("J3").Value = 0.8564255319148937
$ws.Range("K3").Value = 0.4105
$ws.Range("L3").Value = 0.5625
$ws.Range("M3").Value = '[1 1 0 1 0 0 1 0 1 1 1 0 1 1 1 1 1 1 1 1 0 0 1 0]'
$ws.Range("N3").Value = '[0 1 1 1 1 0 0 0 1 1 0 0 1 0 1 1 1 1 0 0 1 0 0 1]'

$ws.Range("A4").Value = 'Pipeline(steps=[(''scaler'', RobustScaler()),
                (''selector'', RandomUnderSampler(random_state=42)),
                (''model'',
                 AdaBoostClassifier(estimator=DecisionTreeClassifier(class_weight=''balanced'',
                                                                     criterion=''entropy'',
                                                                     max_depth=5,
                                                                     max_features=''log2'',
                                                                     min_samples_leaf=3,
                                                                     min_samples_split=4,
                                                                     random_state=42),
                                    n_estimators=5, random_state=42))])'
$ws.Range("B4").Value = 0.6095959595959596
$ws.Range("C4").Value = '{''selector'': RandomUnderSampler(random_state=42), ''scaler'': RobustScaler(), ''model__n_estimators'': 5, ''model__estimator__min_samples_split'': 4, ''model__estimator__min_samples_leaf'': 3, ''model__estimator__max_features'': ''log2'', ''model__estimator__max_depth'': 5, ''model__estimator__criterion'': ''entropy'', ''model__estimator__class_weight'': ''balanced''}'
$ws.Range("D4").Value = 0.9430558065748097
$ws.Range("E4").Value = 0.5020255855255855
$ws.Range("F4").Value = 0.5625
$ws.Range("G4").Value = 0.9872708877627798
$ws.Range("H4").Value = 0.5376944444444445
$ws.Range("I4").Value = 0.6923076923076923
$ws.Range("J4").Value = 0.9036888888888889
$ws.Range("K4").Value = 0.4994000000000001
$ws.Range("L4").Value = 0.4736842105263158
$ws.Range("M4").Value = '[0 1 0 0 1 1 1 1 1 1 1 0 1 1 1 1 1 1 1 1 0 1 1 1]'
$ws.Range("N4").Value = '[0 0 1 1 0 1 0 1 0 0 1 1 0 0 1 1 1 1 1 0 1 1 0 0]'

$ws.Range("A5").Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()),
                (''selector'', RandomUnderSampler(random_state=42)),
                (''model'',
                 AdaBoostClassifier(estimator=DecisionTreeClassifier(criterion=''entropy'',
                                                                     max_depth=5,
                                                                     min_samples_split=3,
                                                                     random_state=42),
                                    random_state=42))])'
$ws.Range("B5").Value = 0.61993006993007
$ws.Range("C5").Value = '{''selector'': RandomUnderSampler(random_state=42), ''scaler'': MinMaxScaler(), ''model__n_estimators'': 50, ''model__estimator__min_samples_split'': 3, ''model__estimator__min_samples_leaf'': 1, ''model__estimator__max_features'': None, ''model__estimator__max_depth'': 5, ''model__estimator__criterion'': ''entropy'', ''model__estimator__class_weight'': None}'
$ws.Range("D5").Value = 0.888061214508038
$ws.Range("E5").Value = 0.4877457875457875
$ws.Range("F5").Value = 0.5925925925925927
$ws.Range("G5").Value = 0.9808044650995831
$ws.Range("H5").Value = 0.5589289682539682
$ws.Range("I5").Value = 0.6153846153846154
$ws.Range("J5").Value = 0.8133877551020408
$ws.Range("K5").Value = 0.4581666666666667
$ws.Range("L5").Value = 0.5714285714285714
$ws.Range("M5").Value = '[0 1 1 0 0 1 0 0 0 0 1 1 1 0 0 1 1 0 1 1 1 1 1 1]'
$ws.Range("N5").Value = '[1 1 1 1 1 0 0 1 0 0 1 0 0 0 1 1 0 0 0 1 0 1 1 1]'

$ws.Range("A6").Value = 'Pipeline(steps=[(''scaler'', RobustScaler()),
                (''selector'', RandomUnderSampler(random_state=42)),
                (''model'',
                 AdaBoostClassifier(estimator=DecisionTreeClassifier(class_weight=''balanced'',
                                                                     max_depth=5,
                                                                     max_features=''sqrt'',
                                                                     min_samples_leaf=6,
                                                                     min_samples_split=6,
                                                                     random_state=42),
                                    n_estimators=10, random_state=42))])'
$ws.Range("B6").Value = 0.6350793650793651
$ws.Range("C6").Value = '{''selector'': RandomUnderSampler(random_state=42), ''scaler'': RobustScaler(), ''model__n_estimators'': 10, ''model__estimator__min_samples_split'': 6, ''model__estimator__min_samples_leaf'': 6, ''model__estimator__max_features'': ''sqrt'', ''model__estimator__max_depth'': 5, ''model__estimator__criterion'': ''gini'', ''model__estimator__class_weight'': ''balanced''}'
$ws.Range("D6").Value = 0.8699615838084128
$ws.Range("E6").Value = 0.5203073482073483
$ws.Range("F6").Value = 0.3636363636363636
$ws.Range("G6").Value = 0.9834527463483009
$ws.Range("H6").Value = 0.6004416666666667
$ws.Range("I6").Value = 0.3636363636363636
$ws.Range("J6").Value = 0.7821730769230767
$ws.Range("K6").Value = 0.4821666666666667
$ws.Range("L6").Value = 0.3636363636363636
$ws.Range("M6").Value = '[1 0 1 1 0 0 0 0 1 0 1 1 0 1 1 0 1 0 0 0 0 0 1 1]'
$ws.Range("N6").Value = '[1 0 0 0 0 1 0 0 0 1 1 0 1 0 0 1 1 1 1 0 1 0 1 0]'
